# "Generate Report for Archive"
#
# The localization status report is refreshed: rows that were previously
# "Ready for handoff" are now back "In Translation". This text lives in a
# shared string that is referenced from the Overview sheet (columns for
# each locale) as well as from each per-locale detail sheet (Status
# column), so updating the cell values updates every occurrence.
#
# Shrinking the status text also lets Excel's column autosizing give back
# a little width on the columns that show it, so the columns that held the
# status text are narrowed accordingly.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"

# Narrow the now-shorter status columns to their new best-fit width.
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
